$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 holds the MIxS field headers; columns AA (27) .. BZ (78) each carry
# a cell comment with the field's description. We are deleting the
# "culture_collection" field (column AA) again, so every field from AB..BZ
# needs to shift one column to the left, both its header value and its
# comment text, and the trailing column (BZ) disappears entirely.

$firstCol = 27  # AA
$lastCol  = 78  # BZ

# 1) Capture the current comment text for every column in the range, before
#    we disturb anything (EntireColumn.Delete shifts cell VALUES but leaves
#    comments anchored to their original cell, so we must move the text
#    ourselves afterwards).
$commentTexts = @{}
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $commentTexts[$col] = $ws.Cells.Item(15, $col).Comment.Text()
}

# 2) Delete the culture_collection column entirely. This removes the
#    "culture_collection" header cell/shared string and shifts every column
#    after it (AB..BZ) one place to the left (AA..BY), also fixing up the
#    sheet dimensions/row spans and dropping the now-empty trailing column.
$ws.Range("AA1").EntireColumn.Delete()

# 3) Re-point each comment so column N shows the description that used to
#    belong to column N+1 (matching the header shift above), for every
#    column except the very last one, which no longer has data.
for ($col = $firstCol; $col -lt $lastCol; $col++) {
    $ws.Cells.Item(15, $col).Comment.Text($commentTexts[$col + 1]) | Out-Null
}

# 4) The comment that used to live on the last column (BZ) is now orphaned
#    (its text was already moved into BY in the previous step), so remove it.
$ws.Cells.Item(15, $lastCol).Comment.Delete()
